$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EMPLOYEE DTR")

# VL BALANCE (row 30) value: "7.4.0" -> "2.0.0"
$ws.Range("C30").Value = "2.0.0"

# SL BALANCE (row 31) value: "7.4.0" -> "0.0.0"
$ws.Range("C31").Value = "0.0.0"
